# Updated cryptos list on Fri Nov  8 07:23:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.777.31"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.896.59"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("D5").Value = "'197.59"
$ws.Range("E5").Value = "  +5.33%  "
$ws.Range("D6").Value = "'595.99"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'0.199"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").Value = "2.895.21"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").Value = "'0.430"
$ws.Range("E11").Value = "  +16.55%  "
$ws.Range("D12").Value = "'0.160"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'4.86"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "3.428.92"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("D15").Value = "75.651.31"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "'0.0000191"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "'27.21"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "2.901.23"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'12.94"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'8.73"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "'369.79"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'4.27"
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'70.98"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "3.051.82"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").Value = "'4.17"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "'9.58"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'0.0000107"
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("D30").Value = "'0.996"
$ws.Range("D31").Value = "'1.40"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'501.99"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "'7.66"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'165.03"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").Value = "'20.09"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'19.62"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "'0.112"
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "'0.103"
$ws.Range("E40").Value = "  +20.25%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'179.20"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "'0.344"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'4.96"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "'1.64"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").Value = "'40.07"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "'1.18"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").Value = "'2.31"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'0.568"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'3.71"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'0.655"
$ws.Range("E51").Value = "  +6.91%  "

# Reset style on cells where a leading apostrophe was used to force text,
# so no stray number-format style gets attached to the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
